$d = $word.ActiveDocument

# 1) Change the heading run text from "Background Theory" to "Background "
#    (keeps the existing Heading2 / numPr / ind paragraph formatting untouched)
$d.Content.Find.Execute("Background Theory", $true, $false, $false, $false, $false, $true, 1, $false, "Background ", 2)

# 2) Find that heading paragraph (its text now starts with "Background ") and
#    insert a fresh paragraph right after it as an insertion point.
$headingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Background `r") {
        $headingPara = $p
    }
}

$rng = $headingPara.Range
$rng.InsertParagraphAfter()

$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -eq $headingPara.Range.End) {
        $targetPara = $p
    }
}
$targetRng = $targetPara.Range

# 3) Replace that placeholder paragraph's contents with the full block of new
#    paragraphs (body text + "1.2 Problem Statement" heading + its body),
#    each carrying explicit Times New Roman / 12pt (sz 24) run formatting and no
#    paragraph style / numbering, matching the authored content exactly.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">A human being is a social animal and has the natural ability to see, listen, speak and interact with the external environment.  </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Unfortunately, there are some people who do not have the ability to </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>interact by speaking</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">. The deaf and dumb population is a result of the physical disability of hearing and speaking. In the recent years, there has been a rapid increase </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>in the number of hearing impaired</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> and speech disabled victims due to birth defects, oral diseases and accidents.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> When a speech impaired person speaks to a normal person, the normal person finds it difficult to understand and asks the deaf-dumb person to show gestures for his/her needs. Dumb persons have their own language to communicate with us; the only thing is we need a translator in between.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Sign language is used by deaf and mute people and it is a communication skill that uses gestures instead of sound to convey meaning simultaneously combining hand shapes, orientation and movement of the hands, arms or body and facial expressions to express fluidly a speaker’s thoughts. But most of the time normal people find it difficult to understand this sign language. This presents a major roadblock for people in the deaf and dumb communities when they try to engage in interaction with others, especially in their educational, social and professional environments. Therefore, it is necessary to have an advance gesture recognition or sign language detection system to bridge this communication gap.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>The people who cannot speak or have lost their ability to speak in some accident, it becomes difficult for them to convey their message within the society. To overcome this, a project called ‘SMART GLOVE’ has been designed. Giving a voice to the voiceless has been a cause that many have championed throughout history, but it’s safe to say that none of those efforts involved packing a bunch of sensors into a glove. The main objective of this pro</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>ject is to help deaf and dumb people by removing communication barrier so they are not restricted in a small social circle and are able to convey their feelings and emotions whenever they want.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Smart glove is based on the wearable technology. It is basically a device which has some specific wearable sensors with phenomenal temperature stability. All the sensors are fitted on a glove which measures the</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> different analog parameters associated with the movement of fingers and orientation of the hand during any particular gesture. These sensors read those particular analog values and coding is done in the microcontroller according to these values to recognize the corresponding sign language. The goal of this project is to develop a portable communication system having multiple sensors for Sign Language Recognition and to translate these gestures into text and sound.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>1.2 Problem Statement</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Deaf and normal person communication is as same as two different persons from different countries using two different languages for communication without any common language which leads to problem in communication. Sign language is the only communication tool used by deaf people to communicate to each other. However, normal people do not understand sign language and this creates a large communication barrier between deaf people and normal people. In addition, the sign language is also not easy to learn due to its natural differences in sentence structure and grammar. Therefore, there is a need to develop a system which can help in translating the sign language into text and voice in order to ensure the effective communication can easily take place in the community.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$targetRng.InsertXML($xml)
